$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header added to the fruit table (column F)
$ws.Range("F1").Value = "new"

# Leave the selection on the newly added header cell, matching the
# saved workbook's cursor position.
$ws.Range("F1").Select() | Out-Null
